$d = $word.ActiveDocument
$enDash = [char]0x2013

# ---- Edit 1: first paragraph gets trailing spaces + 3 new red runs ----
$p1 = $d.Paragraphs(1).Range
# Exclude trailing paragraph mark so we append inside the paragraph
$p1.MoveEnd(1, -1) | Out-Null
$p1.InsertAfter("  ")
$p1.Collapse(0) | Out-Null

$r1 = $p1.Duplicate
$r1.InsertAfter("(This is a change " + $enDash + " Ve")
$r1.Font.Color = 192
$r1.Collapse(0) | Out-Null

$r2 = $r1.Duplicate
$r2.InsertAfter("rsion for branch alternate")
$r2.Font.Color = 192
$r2.Collapse(0) | Out-Null

$r3 = $r2.Duplicate
$r3.InsertAfter(")")
$r3.Font.Color = 192

# ---- Edit 2: add a new, bare, shaded paragraph at the very end of the body ----
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.Collapse(0) | Out-Null
$endRange.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newRange = $newPara.Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$newRange.InsertXML($xml) | Out-Null

Write-Host "Edits complete"
